$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (id) and C (speaker_variant) for rows 2-13.
# Column D (is_prefered) is cleared for all these rows (no more "x" markers,
# and no more "lev distance" implied grouping).
$data = @(
    @{ Row = 2;  B = "#bartolomeus"; C = "Bartolomeus" },
    @{ Row = 3;  B = "#nicodemus";   C = "Nicodemus" },
    @{ Row = 4;  B = "#krispyn";     C = "Krispyn" },
    @{ Row = 5;  B = "#margo";       C = "Margo" },
    @{ Row = 6;  B = "#kriepyn";     C = "Kriepyn" },
    @{ Row = 7;  B = "#nicodemes";   C = "Nicodemes" },
    @{ Row = 8;  B = "#alaradus";    C = "Alaradus" },
    @{ Row = 9;  B = "#konstant";    C = "Konstant" },
    @{ Row = 10; B = "#batolomeus";  C = "Batolomeus" },
    @{ Row = 11; B = "#alardus";     C = "Alardus" },
    @{ Row = 12; B = "#julia";       C = "Julia" },
    @{ Row = 13; B = "#alradus";     C = "Alradus" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = ""
}
